# Fixed damage property type: fixed property money -> diamond
# Updates columns E (MAXHP) and F (MAXMP) for rows 11-70 to a new
# arithmetic progression (500 + (row-10)*50), and moves the active
# selection on the sheet to G30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 11; $row -le 70; $row++) {
    $value = 500 + ($row - 10) * 50
    $ws.Cells.Item($row, 5).Value = $value   # Column E
    $ws.Cells.Item($row, 6).Value = $value   # Column F
}

$ws.Range("G30").Select()
